# Update the two-digit-divided-by-one-digit practice answers in the
# single table of the document. Each data row (1, 5, 9, 13, 17) holds
# five answers across the five columns; replace each cell's text with
# its new value while leaving run formatting (font/size) untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    @("88÷6=14, 4", "39÷7=5, 4", "41÷4=10, 1", "22÷6=3, 4", "24÷3=8, 0"),
    @("35÷4=8, 3", "75÷6=12, 3", "78÷7=11, 1", "29÷3=9, 2", "83÷2=41, 1"),
    @("93÷4=23, 1", "21÷5=4, 1", "37÷8=4, 5", "73÷6=12, 1", "18÷7=2, 4"),
    @("23÷8=2, 7", "14÷3=4, 2", "85÷8=10, 5", "68÷8=8, 4", "98÷6=16, 2"),
    @("22÷6=3, 4", "58÷7=8, 2", "94÷5=18, 4", "37÷7=5, 2", "79÷7=11, 2")
)

$dataRows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $dataRows.Count; $i++) {
    $row = $dataRows[$i]
    $rowValues = $newValues[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
